$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 150
$ws.Range("F3").Value = 157
$ws.Range("F4").Value = 2094
$ws.Range("F5").Value = 4238
$ws.Range("F6").Value = 554
$ws.Range("F7").Value = 1057
$ws.Range("F8").Value = 1320
$ws.Range("F9").Value = 664
$ws.Range("F10").Value = 377
$ws.Range("F11").Value = 2213
$ws.Range("F12").Value = 400
$ws.Range("F13").Value = 662636
$ws.Range("F14").Value = 1644
$ws.Range("F15").Value = 531
$ws.Range("F16").Value = 1466
$ws.Range("F18").Value = 544
$ws.Range("F19").Value = 1281
$ws.Range("F20").Value = 2254
$ws.Range("F21").Value = 1147
$ws.Range("F22").Value = 2706
$ws.Range("F23").Value = 1563
$ws.Range("F24").Value = 828
$ws.Range("F25").Value = 1547
$ws.Range("F26").Value = 532
$ws.Range("F27").Value = 1089
$ws.Range("F28").Value = 286
$ws.Range("F29").Value = 1089
$ws.Range("F31").Value = 83
$ws.Range("F32").Value = 2029
$ws.Range("F33").Value = 1385
$ws.Range("F34").Value = 577
$ws.Range("F35").Value = 1295
$ws.Range("F36").Value = 2622
$ws.Range("F37").Value = 5
$ws.Range("F38").Value = 1148
$ws.Range("F39").Value = 33
$ws.Range("F40").Value = 199
$ws.Range("F41").Value = 2594
$ws.Range("F42").Value = 211
$ws.Range("F43").Value = 990
$ws.Range("F44").Value = 3148
$ws.Range("F45").Value = 1010
$ws.Range("F46").Value = 27
$ws.Range("F47").Value = 881
$ws.Range("F50").Value = 11

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 66
$ws.Range("F6").Value = 471
$ws.Range("F10").Value = 484
$ws.Range("F11").Value = 144685
$ws.Range("F12").Value = 144685
$ws.Range("F18").Value = 231
$ws.Range("F19").Value = 334
$ws.Range("F22").Value = 140
$ws.Range("F26").Value = 579
$ws.Range("F30").Value = 60
$ws.Range("F31").Value = 355
$ws.Range("F32").Value = 275
$ws.Range("F34").Value = 54
$ws.Range("F35").Value = 54
$ws.Range("F38").Value = 208

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 3134
$ws.Range("F5").Value = 247
$ws.Range("F8").Value = 1195
$ws.Range("F10").Value = 1604
$ws.Range("F12").Value = 105
$ws.Range("F13").Value = 1924

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1195
$ws.Range("F5").Value = 150
$ws.Range("F6").Value = 1604
$ws.Range("F8").Value = 157
$ws.Range("F9").Value = 2094
$ws.Range("F10").Value = 105
$ws.Range("F11").Value = 1924
$ws.Range("F12").Value = 4238
$ws.Range("F13").Value = 554
$ws.Range("F14").Value = 1320
$ws.Range("F15").Value = 664
$ws.Range("F16").Value = 377
$ws.Range("F17").Value = 2213
$ws.Range("F19").Value = 662642
$ws.Range("F21").Value = 484
$ws.Range("F22").Value = 1644
$ws.Range("F23").Value = 144685
$ws.Range("F24").Value = 1466
$ws.Range("F26").Value = 544
$ws.Range("F27").Value = 1281
$ws.Range("F28").Value = 2254
$ws.Range("F29").Value = 1147
$ws.Range("F30").Value = 2706
$ws.Range("F31").Value = 1563
$ws.Range("F32").Value = 828
$ws.Range("F34").Value = 1547
$ws.Range("F35").Value = 532
$ws.Range("F36").Value = 140
$ws.Range("F37").Value = 1089
$ws.Range("F38").Value = 1089
$ws.Range("F39").Value = 1385
$ws.Range("F40").Value = 1295
$ws.Range("F41").Value = 2622
$ws.Range("F42").Value = 5
$ws.Range("F43").Value = 1148
$ws.Range("F44").Value = 355
$ws.Range("F45").Value = 275
$ws.Range("F46").Value = 54
$ws.Range("F47").Value = 2594
$ws.Range("F48").Value = 3148
$ws.Range("F49").Value = 208
$ws.Range("F50").Value = 1010
$ws.Range("F53").Value = 11
